$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update item code / name values (A2/P2 share "ST012-Asd" -> "WE002-RELOJ")
$ws.Range("A2").Value = "WE002-RELOJ"
$ws.Range("P2").Value = "WE002-RELOJ"

# Update short code values (B2/T2 share "ST012" -> "WE002")
$ws.Range("B2").Value = "WE002"
$ws.Range("T2").Value = "WE002"

# Update quantity/price value in G2 from 1 to 45
$ws.Range("G2").Value = 45
